# Natmi following Dr Hou advice
# Recompute ligand/receptor-expressing-cell counts (E, K columns) from 1 -> 3
# and the downstream derived statistics that depend on them, for data rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> column letter -> new value
$updates = @{
    2 = @{
        E = 3
        G = 4.523724666666666
        H = 13.571174
        I = 0.4806607624766543
        J = 0.4806607624766543
        K = 3
        M = 4.806204333333334
        N = 14.418613
        O = 0.7287437301541012
        P = 0.7287437301541012
        Q = 21.74194509574022
        R = 195.677505861662
        S = 0.3502785169859515
        T = 0.3502785169859515
    }
    3 = @{
        E = 3
        G = 4.523724666666666
        H = 13.571174
        I = 0.4806607624766543
        J = 0.4806607624766543
        K = 3
        M = 1.788987
        N = 5.366961
        O = 0.2712562698458988
        P = 0.2712562698458988
        Q = 8.092884620245998
        R = 72.835961582214
        S = 0.1303822454907028
        T = 0.1303822454907028
    }
    4 = @{
        E = 3
        G = 2.109481
        H = 6.328443
        I = 0.2241393587371326
        J = 0.2241393587371326
        K = 3
        M = 4.806204333333334
        N = 14.418613
        O = 0.7287437301541012
        P = 0.7287437301541012
        Q = 10.13859672328434
        R = 91.247370509559
        S = 0.1633401523604462
        T = 0.1633401523604462
    }
    5 = @{
        E = 3
        G = 2.109481
        H = 6.328443
        I = 0.2241393587371326
        J = 0.2241393587371326
        K = 3
        M = 1.788987
        N = 5.366961
        O = 0.2712562698458988
        P = 0.2712562698458988
        Q = 3.773834085747
        R = 33.964506771723
        S = 0.06079920637668636
        T = 0.06079920637668636
    }
    6 = @{
        E = 3
        G = 2.778265
        H = 8.334795
        I = 0.2951998787862131
        J = 0.2951998787862131
        K = 3
        M = 4.806204333333334
        N = 14.418613
        O = 0.7287437301541012
        P = 0.7287437301541012
        Q = 13.35290928214833
        R = 120.176183539335
        S = 0.2151250608077034
        T = 0.2151250608077034
    }
    7 = @{
        E = 3
        G = 2.778265
        H = 8.334795
        I = 0.2951998787862131
        J = 0.2951998787862131
        K = 3
        M = 1.788987
        N = 5.366961
        O = 0.2712562698458988
        P = 0.2712562698458988
        Q = 4.970279967554999
        R = 44.732519707995
        S = 0.08007481797850964
        T = 0.08007481797850964
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
